# guide41_dashboad.xlsx — "Add files via upload" commit
#
# Semantic changes reconstructed from the OOXML diff:
#   1. On sheet "p1" (ダッシュボードの構成) and sheet "p2" (言語を切り替える方法),
#      a new blank row is inserted right after the page's title row (row 4),
#      pushing all following content down by one row.
#   2. The two heading strings used on "p2" change from an (invalid) <h4>...</h3>
#      pairing to a consistent <h3>...</h3> pairing:
#        <h4><a name="language menu"></a>言語メニューで選択</h3>
#          -> <h3><a name="language menu"></a>言語メニューで選択</h3>
#        <h4><a name="preferred language"></a>優先言語の設定</h3>
#          -> <h3><a name="preferred language"></a>優先言語の設定</h3>
#   3. The active/selected sheet moves from "p1" to "p2", and each sheet's
#      remembered cell selection changes (p1 -> B10, p2 -> B15).

$wb = $excel.ActiveWorkbook

$p1 = $wb.Worksheets.Item("p1")
$p2 = $wb.Worksheets.Item("p2")

# --- p1: insert a blank row after the title row (row 4) ---------------
$null = $p1.Rows.Item(5).Insert()

# --- p2: insert a blank row after the title row (row 4) ---------------
$null = $p2.Rows.Item(5).Insert()

# --- p2: fix up the two section headings, <h4>...</h3> -> <h3>...</h3> -
$p2.Range("B7").Value = '<h3><a name="language menu"></a>言語メニューで選択</h3>'
$p2.Range("B10").Value = '<h3><a name="preferred language"></a>優先言語の設定</h3>'

# --- selection / active-tab bookkeeping --------------------------------
$null = $p1.Range("B10").Select()
$null = $p2.Range("B15").Select()

# p2 becomes the active (visible/selected) sheet in the saved workbook
$null = $p2.Activate()
